# "add left right wall"
#
# Translation sheet, rows 18-21 (before the edit):
#   18: SingleUseId23 | Default | Left   | LTR | >
#   19: SingleUseId24 | Default | Left   | LTR | <
#   20: SingleUseId25 | Default | Center | LTR | >
#   21: SingleUseId26 | Default | Left   | LTR | New Text
#
# After the edit the old "left wall"/"right wall" rows (18 & 19) are removed
# (rows 20-21 shift up to become 18-19), and a new Score text row is appended:
#   18: SingleUseId25 | Default | Center | LTR | >
#   19: SingleUseId26 | Default | Left   | LTR | New Text
#   20: SingleUseId28 | Default | Left   | LTR | Score<10>
#   21-23: blank row stubs left behind by the row shuffle

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Remove the two "wall" rows; everything below shifts up by two rows.
$ws.Rows("18:19").Delete()

# New trailing row with the Score text id/value.
$ws.Range("B20").Value = "SingleUseId28"
$ws.Range("C20").Value = "Default"
$ws.Range("D20").Value = "Left"
$ws.Range("E20").Value = "LTR"
$ws.Range("F20").Value = "Score<10>"

# Leave behind empty row stubs for 21-23 (matches the trailing blank rows
# introduced by the edit) without stamping any row-height/format override.
$ws.Rows("21:23").Group()
$ws.Rows("21:23").Ungroup()
